$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Prabowo Beri Petunjuk Sangat Detail soal MBG, Terutama Masalah Kebersihan", "2025-09-29T09:05:05+07:00", "Tidak Diketahui", "https://www.liputan6.com/news/read/6171092/prabowo-beri-petunjuk-sangat-detail-soal-mbg-terutama-masalah-kebersihan", "mbg, prabowo"),
    @("Unggah Foto Bareng Donald Trump, Prabowo Yakin Perdamaian Dunia Segera Tercapai", "2025-09-29T07:49:22+07:00", "Tidak Diketahui", "https://www.liputan6.com/news/read/6171054/unggah-foto-bareng-donald-trump-prabowo-yakin-perdamaian-dunia-segera-tercapai", "prabowo"),
    @("Prabowo Rapat Bareng Menteri di Kertanegara Minggu Malam, Bahas MBG hingga Migas", "2025-09-29T01:01:30+07:00", "Tidak Diketahui", "https://www.liputan6.com/news/read/6170969/prabowo-rapat-bareng-menteri-di-kertanegara-minggu-malam-bahas-mbg-hingga-migas", "mbg, prabowo"),
    @("Prabowo Selesaikan Lawatan 4 Negara, Bawa Pulang Investasi Rp 380 Triliun hingga Artefak B...", "2025-09-29T13:15:46+07:00", "Tidak Diketahui", "https://www.liputan6.com/bisnis/read/6171276/prabowo-selesaikan-lawatan-4-negara-bawa-pulang-investasi-rp-380-triliun-hingga-artefak-bersejarah", "prabowo"),
    @("Baru 3 Minggu Jadi Menkeu, Ini Deretan Gebrakan Purbaya", "2025-09-29T12:32:46+07:00", "Tidak Diketahui", "https://www.liputan6.com/bisnis/read/6171279/baru-3-minggu-jadi-menkeu-ini-deretan-gebrakan-purbaya", "purbaya"),
    @("Menkeu Purbaya Tak Naikan Tarif Cukai Rokok, Siapa yang Untung?", "2025-09-29T12:00:15+07:00", "Tidak Diketahui", "https://www.liputan6.com/bisnis/read/6171204/menkeu-purbaya-tak-naikan-tarif-cukai-rokok-siapa-yang-untung", "purbaya"),
    @("Menkeu Purbaya Sidak ke Kantor Pusat BNI, Ada apa?", "2025-09-29T11:40:20+07:00", "Tidak Diketahui", "https://www.liputan6.com/bisnis/read/6171242/menkeu-purbaya-sidak-ke-kantor-pusat-bni-ada-apa", "purbaya")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
    $ws.Cells.Item($row, 5).Value = $data[$i][4]
}
